$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C11").Value = -12.637
$ws.Range("B12").Value = 4.935
$ws.Range("C23").Value = -12.748
$ws.Range("C28").Value = -12.809
$ws.Range("B32").Value = 6.455
$ws.Range("C32").Value = -12.93
$ws.Range("C34").Value = -11.573
$ws.Range("B36").Value = 8.705
$ws.Range("B38").Value = 5.743
$ws.Range("C42").Value = -12.221
$ws.Range("B46").Value = 6.248
$ws.Range("B54").Value = 5.401
$ws.Range("C54").Value = -13.017
$ws.Range("B55").Value = 4.763
$ws.Range("B67").Value = 5.505
$ws.Range("B69").Value = 5.367
$ws.Range("B72").Value = 5.697
$ws.Range("B91").Value = 5.296
$ws.Range("C97").Value = -11.335
$ws.Range("B99").Value = 5.11
$ws.Range("C99").Value = -12.225
$ws.Range("C101").Value = -12.17
$ws.Range("B104").Value = 8.358000000000001
